{"js": "// Fill in the first empty paragraph right after the \"Write Up\" title with\n// the two new intro paragraphs, followed by a new Heading 1 paragraph\n// announcing the article title. The two trailing empty paragraphs are left\n// untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The first paragraph is the \"Write Up\" title; the very next paragraph is\n// the empty one that gets filled in with the new content.\nconst target = paragraphs.items[1];\n\ntarget.insertText(\n  \"This week, we will be looking at the Poly Build tool, which we can get to by going to the Tool box in Edit mode while working in Blender. This tool is great for creating new topology from scratch, and combines a few tools together into one, in order to work more efficiently in Blender.\",\n  Word.InsertLocation.replace\n);\n\nconst introParagraph = target.insertParagraph(\n  \"So, if this sounds at all interesting to you then please join us for our brand-new article this week entitled:\",\n  Word.InsertLocation.after\n);\n\nconst headingParagraph = introParagraph.insertParagraph(\n  \"The Poly Build Tool\",\n  Word.InsertLocation.after\n);\nheadingParagraph.style = \"Heading 1\";\n\nawait context.sync();\n", "ps1": "# Fill in the first empty paragraph right after the \"Write Up\" title with\n# the two new intro paragraphs, followed by a new Heading 1 paragraph\n# announcing the article title. The two trailing empty paragraphs are left\n# untouched.\n\n$d = $word.ActiveDocument\n\n# Paragraph 1 is the \"Write Up\" title; paragraph 2 is the empty paragraph\n# right after it that gets filled in with the new content.\n$target = $d.Paragraphs.Item(2)\n$target.Range.Text = \"This week, we will be looking at the Poly Build tool, which we can get to by going to the Tool box in Edit mode while working in Blender. This tool is great for creating new topology from scratch, and combines a few tools together into one, in order to work more efficiently in Blender.\"\n$target.Range.InsertParagraphAfter()\n\n$introParagraph = $d.Paragraphs.Item(3)\n$introParagraph.Range.Text = \"So, if this sounds at all interesting to you then please join us for our brand-new article this week entitled:\"\n$introParagraph.Range.InsertParagraphAfter()\n\n$headingParagraph = $d.Paragraphs.Item(4)\n$headingParagraph.Range.Text = \"The Poly Build Tool\"\n$headingParagraph.Style = \"Heading 1\"\n"}
